$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''25.973.75'
$ws.Range('D3').Value = '''1.636.96'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').Value = '''214.56'
$ws.Range('D6').Value = '''0.5088'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').Value = '''0.2561'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').Value = '''0.06345'
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('D10').Value = '''19.60'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '''4.267'
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').Value = '''1.626.94'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').Value = '''0.5428'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '''0.0₅7699'
$ws.Range('E15').Value = '  -2.38%  '
$ws.Range('D16').Value = '''63.88'
$ws.Range('E16').Value = '  -1.69%  '
$ws.Range('D17').Value = '''25.985.19'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '''1.001'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('D19').Value = '''198.99'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').Value = '''4.412'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '''9.896'
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').Value = '''6.047'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').Value = '''1.003'
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('D24').Value = '''1.889'
$ws.Range('E24').Value = '  +0.95%  '
$ws.Range('D25').Value = '''141.29'
$ws.Range('E25').Value = '  +0.21%  '
$ws.Range('E26').Value = '  +4.95%  '
$ws.Range('D27').Value = '''6.822'
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('D28').Value = '''15.56'
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('D29').Value = '''1.233'
$ws.Range('E29').Value = '  -0.93%  '
$ws.Range('D30').Value = '''0.04896'
$ws.Range('E30').Value = '  -2.15%  '
$ws.Range('D31').Value = '''3.253'
$ws.Range('D32').Value = '''3.169'
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('D33').Value = '''1.527'
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('D34').Value = '''2.367'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').Value = '''0.9077'
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('D36').Value = '''2.582'
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('D37').Value = '''1.127.41'
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('D38').Value = '''0.5452'
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('D39').Value = '''0.01563'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').Value = '''1.001'
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('D41').Value = '''2.524'
$ws.Range('E41').Value = '  -1.44%  '
$ws.Range('D42').Value = '''0.8108'
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('E43').Value = '  +2.29%  '
$ws.Range('D44').Value = '''99.02'
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('D45').Value = '''5.414'
$ws.Range('E45').Value = '  -5.01%  '
$ws.Range('D46').Value = '''1.774.74'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('E48').Value = '  -0.17%  '
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').Value = '''0.05105'
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('E51').Value = '  -0.27%  '
